$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C12").Value = 5
$ws.Range("D12").Value = 5
$ws.Range("E12").Value = 5
$ws.Range("F12").Value = 5

$ws.Range("G12").Select()
